$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had an accidental duplicate entry at row 138
# (id_joueur=14, id_mecanique=47 -- already present a few rows earlier/later);
# remove it so every row below shifts up by one (B301 -> B300 data range).
$ws.Rows("138:138").Delete()

# Column D only ever held a leftover width/formatting definition from an old
# "Dico" helper column -- no real data lived in it. Deleting it removes that
# stale column and breaks the "Dico" defined name (it used to point at
# Feuil1!$D:$D), which Excel turns into a #REF! error once its target
# column is gone.
$ws.Columns("D:D").Delete()
$wb.Names.Item("Dico").RefersTo = "=Feuil1!#REF!"

# Leave the selection where the user was last working.
$ws.Range("H147").Select() | Out-Null
